$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties columns (AD1:AF1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, centered) from an existing header cell (AC1) to the new ones
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-40: team record is the same for every player row
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 101  # column AD
    $ws.Cells.Item($r, 31).Value = 61   # column AE
    $ws.Cells.Item($r, 32).Value = 0    # column AF
}
